# Update Junction Flooding data table: replace rows 2-5 with new 1000-sample
# readings, drop the old row 6, and widen several data columns from 7 to 8
# characters (custom accuracy formatting pass).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sensor readings for rows 2-5 (column 1 = A ... column 34 = AH)
$data = @{
    2 = @(45075.50694444445, 12.77, 8.425000000000001, 3.488, 27.776, 20.765, 9.82, 29.191, 15.707, 6.218, 9.112, 10.928, 11.76, 3.255, 10.151, 13.871, 9.122999999999999, 2.77, 1.592, 147.125, 28.067, 9.369999999999999, 18.004, 9.226000000000001, 2.728, 15.858, 8.276, 7.652, 9.069000000000001, 11.378, 2.869, 26.432, 4.957, 11.713)
    3 = @(45075.51388888889, 18.102, 13.159, 1.817, 39.655, 31.61, 14.122, 53.101, 22.106, 9.574, 14.003, 15.847, 16.923, 4.589, 14.287, 20.129, 12.349, 1.395, 1.064, 210.134, 39.929, 13.187, 26.485, 13.819, 2.463, 26.538, 11.648, 10.481, 12.341, 16.576, 1.165, 48.637, 7.277, 16.487)
    4 = @(45075.52083333334, 23.416, 17.348, 1.539, 51.193, 41.514, 18.347, 71.07599999999999, 28.505, 12.583, 18.542, 20.512, 21.831, 5.918, 18.422, 26.151, 15.635, 1.022, 1.02, 273.095, 51.516, 17.004, 34.51, 18.071, 2.806, 34.824, 15.02, 13.365, 15.723, 21.515, 0.733, 64.77, 9.52, 21.259)
    5 = @(45075.52777777778, 2.78, 1.93, 0.64, 6.26, 4.56, 2.13, 16.64, 3.49, 1.51, 1.89, 2.5, 2.75, 0.75, 2.26, 3.24, 2.16, 0.64, 0.32, 27.12, 6.8, 2.08, 4.43, 2.16, 0.59, 7.43, 1.84, 1.77, 2.06, 2.59, 0.54, 15.88, 1.05, 2.61)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

# Drop the old trailing row (row 6) entirely; dimension shrinks to AH5
$ws.Rows(6).Delete()

# Widen these data columns from 7 to 8 characters
$wideCols = @(2, 3, 7, 10, 11, 12, 13, 15, 16, 17, 22, 24, 27, 28, 29, 30, 34)
foreach ($c in $wideCols) {
    $ws.Cells.Item(1, $c).EntireColumn.ColumnWidth = 7.14
}
